$geneMap = @{
    "Irf6" = "ENSMUSG00000026638"
    "Bcl11a" = "ENSMUSG00000000861"
    "Hey1" = "ENSMUSG00000040289"
    "Foxo4" = "ENSMUSG00000042903"
    "Bcl11b" = "ENSMUSG00000048251"
    "E2f1" = "ENSMUSG00000027490"
    "Zfp362" = "ENSMUSG00000028799"
    "Sox11" = "ENSMUSG00000063632"
    "Foxp4" = "ENSMUSG00000023991"
    "Zfp933" = "ENSMUSG00000059423"
    "Zbtb7b" = "ENSMUSG00000028042"
    "Zfp367" = "ENSMUSG00000044934"
    "Hoxb5" = "ENSMUSG00000038700"
    "Nr6a1" = "ENSMUSG00000063972"
    "Elk1" = "ENSMUSG00000009406"
    "Hsf4" = "ENSMUSG00000033249"
    "Hoxc4" = "ENSMUSG00000075394"
    "Nfix" = "ENSMUSG00000001911"
    "Six1" = "ENSMUSG00000051367"
    "Ahr" = "ENSMUSG00000019256"
    "Yy1" = "ENSMUSG00000021264"
    "Hoxc8" = "ENSMUSG00000001657"
    "Rfx5" = "ENSMUSG00000005774"
    "Srebf2" = "ENSMUSG00000022463"
    "Maff" = "ENSMUSG00000042622"
    "Pou6f1" = "ENSMUSG00000009739"
    "Mypop" = "ENSMUSG00000048481"
    "Tcf12" = "ENSMUSG00000032228"
    "Csrnp3" = "ENSMUSG00000044647"
    "Atf3" = "ENSMUSG00000026628"
    "Twist2" = "ENSMUSG00000007805"
    "Nfyc" = "ENSMUSG00000032897"
    "Etv3" = "ENSMUSG00000003382"
    "Tfdp2" = "ENSMUSG00000032411"
    "Aff3" = "ENSMUSG00000037138"
    "Atf5" = "ENSMUSG00000038539"
    "Tcerg1" = "ENSMUSG00000024498"
    "Smad3" = "ENSMUSG00000032402"
    "Stat3" = "ENSMUSG00000004040"
    "Hoxa2" = "ENSMUSG00000014704"
    "Foxp3" = "ENSMUSG00000039521"
    "Zbtb16" = "ENSMUSG00000066687"
    "Zfp366" = "ENSMUSG00000050919"
    "Gata1" = "ENSMUSG00000031162"
    "Hoxd9" = "ENSMUSG00000043342"
    "Ecsit" = "ENSMUSG00000066839"
    "Fosl1" = "ENSMUSG00000024912"
    "Hinfp" = "ENSMUSG00000032119"
    "Zfp738" = "ENSMUSG00000048280"
    "Pbx4" = "ENSMUSG00000031860"
    "Rora" = "ENSMUSG00000032238"
    "Nfil3" = "ENSMUSG00000056749"
    "Pax8" = "ENSMUSG00000026976"
    "Zfp748" = "ENSMUSG00000095432"
    "Relb" = "ENSMUSG00000002983"
    "Onecut2" = "ENSMUSG00000045991"
    "Nfe2l3" = "ENSMUSG00000029832"
    "Mef2d" = "ENSMUSG00000001419"
    "E2f5" = "ENSMUSG00000027552"
    "Lmx1b" = "ENSMUSG00000038765"
    "Nr1i2" = "ENSMUSG00000022809"
    "Elk3" = "ENSMUSG00000008398"
    "Glis3" = "ENSMUSG00000052942"
    "Spic" = "ENSMUSG00000004359"
    "Foxc2" = "ENSMUSG00000046714"
    "Batf" = "ENSMUSG00000034266"
    "Hoxb6" = "ENSMUSG00000000690"
    "Sox9" = "ENSMUSG00000000567"
    "Hlf" = "ENSMUSG00000003949"
    "Hif3a" = "ENSMUSG00000004328"
    "Tcf7" = "ENSMUSG00000000782"
    "Foxj1" = "ENSMUSG00000034227"
    "Hoxd4" = "ENSMUSG00000101174"
    "Cebpb" = "ENSMUSG00000056501"
    "Erf" = "ENSMUSG00000040857"
    "Myb" = "ENSMUSG00000019982"
    "Nrl" = "ENSMUSG00000040632"
    "Zfp708" = "ENSMUSG00000058883"
    "Tbx21" = "ENSMUSG00000001444"
    "Meox1" = "ENSMUSG00000001493"
    "Rfx4" = "ENSMUSG00000020037"
    "Klf12" = "ENSMUSG00000072294"
    "Nfe2" = "ENSMUSG00000058794"
    "Neurog3" = "ENSMUSG00000044312"
    "Zfp408" = "ENSMUSG00000075040"
    "Creb3l4" = "ENSMUSG00000027938"
    "Spi1" = "ENSMUSG00000002111"
    "Foxo6" = "ENSMUSG00000052135"
    "Irf8" = "ENSMUSG00000041515"
}

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count

    $ws.Cells.Item(1, 7).Value = "gene_id"

    for ($r = 2; $r -le $rowCount; $r++) {
        $symbol = $ws.Cells.Item($r, 6).Value2
        if ($geneMap.ContainsKey($symbol)) {
            $ws.Cells.Item($r, 7).Value = $geneMap[$symbol]
        }
    }
}
